$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-11, columns A (datetime serial), B (humidity),
# C (year), D (month), E (day), F (hour), G (minute), H (dayofyear), I (dayofweek)
$data = @(
    @(45701.35902777778, 100,             2025, 2, 13,  8, 37, 44, 3),
    @(45702.72986111111, 96.31,           2025, 2, 14, 17, 31, 45, 4),
    @(45711.00486111111, 96.98,           2025, 2, 23,  0,  7, 54, 6),
    @(45711.6625,        100,             2025, 2, 23, 15, 54, 54, 6),
    @(45711.96180555555, 100,             2025, 2, 23, 23,  5, 54, 6),
    @(45713.64097222222, 99.63,           2025, 2, 25, 15, 23, 56, 1),
    @(45714.91527777778, 96.83,           2025, 2, 26, 21, 58, 57, 2),
    @(45721.53680555556, 98.31,           2025, 3,  5, 12, 53, 64, 2),
    @(45721.56458333333, 98.81,           2025, 3,  5, 13, 33, 64, 2),
    @(45721.89652777778, 99.26000000000001, 2025, 3, 5, 21, 31, 64, 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
}
